# Update the "2859TH daily flow" date range to reflect the latest data
# (notes after data meeting with kate)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = "1/1/2000 to present"

# Move the active selection to D6 (where the edit was made)
$ws.Range("D6").Select()
